# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the rows whose quoted figures changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.994.42"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.633.23"
$ws.Range("E3").Value = "  +3.79%  "
$ws.Range("D5").Value = "517.96"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "144.51"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "2.658.39"
$ws.Range("E9").Value = "  +4.59%  "
$ws.Range("D10").Value = "6.25"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "3.100.86"
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("D15").Value = "58.932.48"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "20.88"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "2.657.39"
$ws.Range("E18").Value = "  +4.76%  "
$ws.Range("D19").Value = "347.40"
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "61.74"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "0.0₃0802"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "7.11"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  +7.29%  "
$ws.Range("D32").Value = "18.94"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").Value = "149.57"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "0.974"
$ws.Range("E35").Value = "  +6.44%  "
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").Value = "1.14"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("D38").Value = "36.74"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").Value = "0.841"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("E40").Value = "  +4.92%  "
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").Value = "278.18"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "0.611"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "19.64"
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "10.30"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "1.987.75"
$ws.Range("E50").Value = "  +4.48%  "
$ws.Range("D51").Value = "4.67"
$ws.Range("E51").Value = "  +3.35%  "
